$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# URL: http://ibm.com/... -> http://linuxforhealth.org/...
$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/legal-document-system"

# Version: 7.0.0 -> 8.0.0
$ws1.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00 (keep as text)
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$ws1.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Legal Document S" ---
$ws2 = $wb.Worksheets.Item("Include from Legal Document S")

# System URI: http://ibm.com/... -> http://linuxforhealth.org/...
$ws2.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/legal-document-system"
